{"js": "// Apply the \"make bill downloadable and add routes\" billing data edit.\n//\n// The document is a billing summary with several placeholder values\n// (the literal text \"123\") that must be replaced with real figures, plus\n// a handful of unique numeric/text fields that change to their final\n// values. Because the string \"123\" appears multiple times, those hits\n// are resolved by their position in document order; every other value\n// is unique in the document and is replaced directly by searching for it.\n\nconst body = context.document.body;\n\n// --- 1) Unique (unambiguous) text/number replacements -----------------\nconst uniqueReplacements = [\n  [\"Jan Sch\u00f6nfeld\", \"G\u00fcnther Netzer\"],\n  [\"68.68\", \"55.83\"],\n  [\"36.9\", \"6\"],\n  [\"7.38\", \"1.2\"],\n  [\"18.45\", \"2.25\"],\n  [\"131.41\", \"65.28\"],\n  [\"24.97\", \"12.4\"],\n  [\"156.38\", \"77.68\"],\n];\n\nfor (const [findText, replaceText] of uniqueReplacements) {\n  const results = body.search(findText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length !== 1) {\n    throw new Error(\n      `Expected exactly 1 match for \"${findText}\", found ${results.items.length}`\n    );\n  }\n\n  results.items[0].insertText(replaceText, Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// --- 2) The repeated \"123\" placeholders, resolved by document order ---\nconst placeholderResults = body.search(\"123\", { matchCase: true });\nplaceholderResults.load(\"items\");\nawait context.sync();\n\nconst placeholderReplacements = [\n  \"1 XVII XCA 111\", // name suffix / house number\n  \"12 345 678 910\", // Ident. Nr.\n  \"100\", // Betreuungsarbeit minutes\n  \"20\", // Fahrkosten Km\n  \"20\", // Telefongeb. Einh.\n  \"15\", // Kopien x\n  \"10\", // Porto\n];\n\nif (placeholderResults.items.length !== placeholderReplacements.length) {\n  throw new Error(\n    `Expected ${placeholderReplacements.length} \"123\" matches, found ${placeholderResults.items.length}`\n  );\n}\n\nfor (let i = 0; i < placeholderResults.items.length; i++) {\n  placeholderResults.items[i].insertText(\n    placeholderReplacements[i],\n    Word.InsertLocation.replace\n  );\n}\nawait context.sync();\n", "ps1": "# Apply the \"make bill downloadable and add routes\" billing data edit.\n#\n# The document is a billing summary with several placeholder values (the\n# literal text \"123\") that must be replaced with real figures, plus a\n# handful of unique numeric/text fields that change to their final\n# values. Because the string \"123\" appears multiple times, those hits\n# are resolved by their position in document order; every other value\n# is unique in the document and is replaced directly via Find & Replace.\n\n$d = $word.ActiveDocument\n\n# --- 1) Unique (unambiguous) text/number replacements ------------------\n$uniqueReplacements = @(\n  @(\"Jan Sch\u00f6nfeld\", \"G\u00fcnther Netzer\"),\n  @(\"68.68\", \"55.83\"),\n  @(\"36.9\", \"6\"),\n  @(\"7.38\", \"1.2\"),\n  @(\"18.45\", \"2.25\"),\n  @(\"131.41\", \"65.28\"),\n  @(\"24.97\", \"12.4\"),\n  @(\"156.38\", \"77.68\")\n)\n\nforeach ($pair in $uniqueReplacements) {\n  $findText = $pair[0]\n  $replaceText = $pair[1]\n  $rng = $d.Content\n  $found = $rng.Find.Execute($findText, $false, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)\n  if (-not $found) {\n    throw \"Could not find expected unique text: $findText\"\n  }\n}\n\n# --- 2) The repeated \"123\" placeholders, resolved by document order ----\n$placeholderReplacements = @(\n  \"1 XVII XCA 111\",  # name suffix / house number\n  \"12 345 678 910\",  # Ident. Nr.\n  \"100\",             # Betreuungsarbeit minutes\n  \"20\",              # Fahrkosten Km\n  \"20\",              # Telefongeb. Einh.\n  \"15\",              # Kopien x\n  \"10\"               # Porto\n)\n\n$rng = $d.Content\n$idx = 0\nwhile ($idx -lt $placeholderReplacements.Length -and $rng.Find.Execute(\"123\")) {\n  $rng.Text = $placeholderReplacements[$idx]\n  $rng.Collapse(0)  # wdCollapseEnd\n  $idx = $idx + 1\n}\n\nif ($idx -ne $placeholderReplacements.Length) {\n  throw \"Expected $($placeholderReplacements.Length) '123' placeholders, replaced $idx\"\n}\n"}
